# Auto-generated: append 42 survey rows (2656-2697) to 'Sheet 1'
# and add the new shared string 'taianivargad@hotmail.com' used in
# the Email Address column of rows 2674-2675.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstNewRow = 2656
$lastNewRow = 2697

# Column A keeps the same "yyyy/mm/dd hh:mm:ss" date style (s="1") used
# by every existing Timestamp cell above - copy it down before writing
# the new values so no new style entries are created.
$ws.Range("A2655").Copy()
$ws.Range($ws.Cells.Item($firstNewRow, 1), $ws.Cells.Item($lastNewRow, 1)).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$rows = @(
  @(2656, 44001.4587975, "fernanda.stringhi@maededeus.com.br", "Hospital M" + [char]0x00E3 + "e de Deus - UTI ADULTO - TIPO I", 55, 5, 42, 3, 5, "utipoa"),
  @(2657, 44001.4969082986, "lauren.ghion@santacasa.org.br", "Complexo Hospitalar Santa Casa - UTI ADULTO - TIPO III", 87, 0, 69, 1, 7, "utipoa"),
  @(2658, 44001.4977388079, "lauren.ghion@santacasa.org.br", "Complexo Hospitalar Santa Casa - UTI PEDIATRICA - TIPO III", 37, 3, 37, 0, 0, "utipoa"),
  @(2659, 44001.5071293519, "marcosboniatti@gmail.com", "Hospital Cristo Redentor - UTI ADULTO - TIPO III", 29, 0, 26, 0, 0, "utipoa"),
  @(2660, 44001.5600644329, "francojw66@yahoo.com.br", "Hospital de Pronto Socorro de Porto Alegre - UTI DE QUEIMADOS", 4, 0, 4, 0, 0, "utipoa"),
  @(2661, 44001.5605259722, "francojw66@yahoo.com.br", "Hospital de Pronto Socorro de Porto Alegre - UTI DE QUEIMADOS", 4, 0, 4, 0, 0, "utipoa"),
  @(2662, 44001.5878303009, "ccih@hpa.org.br", "Hospital Porto Alegre - UTI ADULTO - TIPO II", 7, 0, 7, 0, 1, "utipoa"),
  @(2663, 44001.6375027894, "joao.krauzer@hmv.org.br", "Hospital Moinhos de Vento - UTI PEDIATRICA - TIPO III", 11, 0, 7, 1, 0, "utipoa"),
  @(2664, 44001.6784255671, "lucirestelatto@gmail.com", "Hospital Divina Provid" + [char]0x00EA + "ncia - UTI ADULTO - TIPO II", 16, 0, 11, 0, 1, "utipoa"),
  @(2665, 44001.6873026157, "pedrocomerlato@hotmail.com", "Hospital Independ" + [char]0x00EA + "ncia - UTI ADULTO - TIPO II", 10, 0, 10, 0, 0, "utipoa"),
  @(2666, 44001.7217681019, "taianivatgas@hotmail.com", "Hospital Nossa Senhora da Concei" + [char]0x00E7 + [char]0x00E3 + "o - UTI ADULTO - TIPO III", 75, 0, 70, 3, 29, "utipoa"),
  @(2667, 44001.7219416088, "taianivatgas@hotmail.com", "Hospital Nossa Senhora da Concei" + [char]0x00E7 + [char]0x00E3 + "o - UTI ADULTO - TIPO III", 75, 0, 70, 3, 29, "utipoa"),
  @(2668, 44001.7233820949, "taianivargas@hotmail.com", "Hospital Nossa Senhora da Concei" + [char]0x00E7 + [char]0x00E3 + "o - UTI ADULTO - TIPO III", 75, 0, 70, 3, 20, "utipoa"),
  @(2669, 44001.7235947222, "taianivargas@hotmail.com", "Hospital Nossa Senhora da Concei" + [char]0x00E7 + [char]0x00E3 + "o - UTI ADULTO - TIPO III", 75, 0, 70, 3, 20, "utipoa"),
  @(2670, 44001.7242786227, "taianivargas@hotmail.com", "Hospital Nossa Senhora da Concei" + [char]0x00E7 + [char]0x00E3 + "o - UTI ADULTO - TIPO III", 75, 0, 70, 3, 20, "utipoa"),
  @(2671, 44001.7243982639, "taianivargas@hotmail.com", "Hospital Nossa Senhora da Concei" + [char]0x00E7 + [char]0x00E3 + "o - UTI ADULTO - TIPO III", 75, 0, 70, 3, 20, "utipoa"),
  @(2672, 44001.7256667014, "taianivargas@hotmail.com", "Hospital Nossa Senhora da Concei" + [char]0x00E7 + [char]0x00E3 + "o - UTI ADULTO - TIPO III", 75, 0, 70, 3, 20, "utipoa"),
  @(2673, 44001.7257756482, "taianivargas@hotmail.com", "Hospital Nossa Senhora da Concei" + [char]0x00E7 + [char]0x00E3 + "o - UTI ADULTO - TIPO III", 75, 0, 70, 3, 20, "utipoa"),
  @(2674, 44001.7349621412, "taianivargad@hotmail.com", "Hospital Nossa Senhora da Concei" + [char]0x00E7 + [char]0x00E3 + "o - UTI ADULTO - TIPO III", 75, 0, 70, 3, 20, "utipoa"),
  @(2675, 44001.735082338, "taianivargad@hotmail.com", "Hospital Nossa Senhora da Concei" + [char]0x00E7 + [char]0x00E3 + "o - UTI ADULTO - TIPO III", 75, 0, 70, 3, 20, "utipoa"),
  @(2676, 44001.7356382176, "andre.machado@hed.com.br", "Hospital Ernesto Dorenelles - UTI ADULTO - TIPO III", 40, 0, 32, 4, 5, "utipoa"),
  @(2677, 44001.7383912153, "renatafarinon@yahoo.com.br", "Hospital Santa Ana - UTI ADULTO - TIPO II", 10, 0, 6, 0, 0, "utipoa"),
  @(2678, 44001.7709613889, "braun.luiz@gmail.com", "Hospital Nossa Senhora da Concei" + [char]0x00E7 + [char]0x00E3 + "o - UTI PEDIATRICA - TIPO II", 18, 6, 9, 3, 0, "utipoa"),
  @(2679, 44001.8073495023, "mscanabarro@terra.com.br", "Hospital Femina - UTI ADULTO - TIPO II", 6, 0, 4, 0, 0, "utipoa"),
  @(2680, 44001.8103806944, "smarcos@ghc.com.br", "Hospital Femina - UTI ADULTO - TIPO II", 6, 0, 3, 0, 0, "utipoa"),
  @(2681, 44001.8465621065, "mscanabarro@terra.com.br", "Hospital Femina - UTI ADULTO - TIPO II", 6, 0, 4, 0, 0, "utipoa"),
  @(2682, 44002.064219213, "joao.krauzer@hmv.org.br", "Hospital Moinhos de Vento - UTI PEDIATRICA - TIPO III", 11, 0, 8, 3, 0, "utipoa"),
  @(2683, 44002.3248381134, "analuizafilipini@gmail.com", "Hospital S" + [char]0x00E3 + "o Lucas - UTI ADULTO - TIPO III", 59, 0, 49, 2, 2, "utipoa"),
  @(2684, 44002.3618189468, "dralubarcellos@gmail.com", "Hospital de Pronto Socorro de Porto Alegre - UTI PEDIATRICA - TIPO III", 8, 0, 5, 0, 0, "utipoa"),
  @(2685, 44002.3751276042, "renatafarinon@yahoo.com.br", "Hospital Santa Ana - UTI ADULTO - TIPO II", 10, 0, 5, 0, 0, "utipoa"),
  @(2686, 44002.3763680093, "andre.machado@hed.com.br", "Hospital Ernesto Dorenelles - UTI ADULTO - TIPO III", 40, 0, 32, 4, 5, "utipoa"),
  @(2687, 44002.3826608102, "smarcos@ghc.com.br", "Hospital Femina - UTI ADULTO - TIPO II", 6, 0, 4, 0, 0, "utipoa"),
  @(2688, 44002.3893542477, "fnagel@hcpa.edu.br", "Hospital de Cl" + [char]0x00ED + "nicas de Porto Alegre - UTI ADULTO - TIPO III", 118, 1, 97, 2, 38, "utipoa"),
  @(2689, 44002.3895087153, "lufacchi@uol.com.br", "Hospital Vila Nova - UTI ADULTO - TIPO II", 20, 0, 19, 1, 1, "utipoa"),
  @(2690, 44002.3904060764, "fnagel@hcpa.edu.br", "Hospital de Cl" + [char]0x00ED + "nicas de Porto Alegre - UTI ADULTO - TIPO III", 118, 1, 98, 9, 35, "utipoa"),
  @(2691, 44002.391707037, "leandra@portoalegre.rs.gov.br", "Hospital Materno Infantil Presidente Vargas - UTI PEDIATRICA - TIPO II", 12, 0, 6, 0, 0, "utipoa"),
  @(2692, 44002.3983591667, "andre.machado@hed.com.br", "Hospital Ernesto Dorenelles - UTI ADULTO - TIPO III", 40, 0, 32, 5, 5, "utipoa"),
  @(2693, 44002.4260890393, "cdalmora@hcpa.edu.br", "Hospital de Cl" + [char]0x00ED + "nicas de Porto Alegre - UTI ADULTO - TIPO III", 13, 0, 12, 0, 0, "utipoa"),
  @(2694, 44002.439876412, "braun.luiz@gmail.com", "Hospital Nossa Senhora da Concei" + [char]0x00E7 + [char]0x00E3 + "o - UTI PEDIATRICA - TIPO II", 18, 6, 9, 3, 0, "utipoa"),
  @(2695, 44002.4491817361, "roseuti@gmail.com", "Hospital Moinhos de Vento - UTI ADULTO - TIPO III", 56, 0, 48, 2, 10, "utipoa"),
  @(2696, 44002.4541492708, "fnagel@hcpa.edu.br", "Hospital de Cl" + [char]0x00ED + "nicas de Porto Alegre - UTI ADULTO - TIPO III", 118, 1, 98, 9, 35, "utipoa"),
  @(2697, 44002.4600214352, "dralubarcellos@gmail.com", "Hospital de Pronto Socorro de Porto Alegre - UTI PEDIATRICA - TIPO III", 8, 0, 6, 0, 0, "utipoa")
)

foreach ($r in $rows) {
  $rowNum = $r[0]
  $ws.Cells.Item($rowNum, 1).Value2 = $r[1]   # A: Timestamp
  $ws.Cells.Item($rowNum, 2).Value  = $r[2]   # B: Email Address
  $ws.Cells.Item($rowNum, 3).Value  = $r[3]   # C: Local Informante
  $ws.Cells.Item($rowNum, 4).Value2 = $r[4]   # D: Quantidade de Leitos
  $ws.Cells.Item($rowNum, 5).Value2 = $r[5]   # E: Leitos bloqueados
  $ws.Cells.Item($rowNum, 6).Value2 = $r[6]   # F: Quantidade de Pacientes
  $ws.Cells.Item($rowNum, 7).Value2 = $r[7]   # G: Pacientes COVID 19 suspeitos
  $ws.Cells.Item($rowNum, 8).Value2 = $r[8]   # H: Pacientes COVID 19 confirmados
  $ws.Cells.Item($rowNum, 9).Value  = $r[9]   # I: Senha de Validacao
}

